$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 already carries the correct per-column cell formatting (same as rows 2-4).
# Copy that formatting down into the new row 5 before populating both rows with data.
# (Row 5 has no K cell, so only copy columns A:J.)
$ws.Range("A6:J6").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New row 5: Crumpet GEF
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 (was blank before): Scone GEF
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Final cursor/selection position left on D7 (matches the saved view state)
$ws.Range("D7").Select() | Out-Null
